$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.84"
$ws.Range("E2").Value = "'0.30%"
$ws.Range("D3").Value = "'38.48"
$ws.Range("E3").Value = "'7.54%"
$ws.Range("D4").Value = "'5.086"
$ws.Range("E4").Value = "'0.81%"
$ws.Range("D5").Value = "'0.08065"
$ws.Range("E5").Value = "'0.38%"
$ws.Range("D6").Value = "'1.943"
$ws.Range("E6").Value = "'4.52%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.186"
$ws.Range("E7").Value = "'1.54%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.934"
$ws.Range("E8").Value = "'2.15%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9293"
$ws.Range("E9").Value = "'0.51%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1459"
$ws.Range("E10").Value = "'14.82%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1931"
$ws.Range("E11").Value = "'1.72%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08954"
$ws.Range("E12").Value = "'-0.59%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03503"
$ws.Range("E13").Value = "'1.98%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09793"
$ws.Range("E14").Value = "'-0.68%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001402"
$ws.Range("E15").Value = "'-0.08%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005892"
$ws.Range("E16").Value = "'-5.84%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.724"
$ws.Range("E17").Value = "'-3.68%"
$ws.Range("D18").Value = "'3.410"
$ws.Range("E18").Value = "'3.05%"
$ws.Range("D19").Value = "'0.3462"
$ws.Range("E19").Value = "'1.41%"
$ws.Range("D20").Value = "'0.1310"
$ws.Range("E20").Value = "'0.80%"
$ws.Range("D21").Value = "'4.783"
$ws.Range("E21").Value = "'-0.58%"
$ws.Range("E22").Value = "'2.92%"
$ws.Range("E23").Value = "'-0.01%"
$ws.Range("D24").Value = "'0.001234"
$ws.Range("E24").Value = "'0.44%"
$ws.Range("D25").Value = "'0.004262"
$ws.Range("E25").Value = "'-12.24%"
$ws.Range("D27").Value = "'0.0001301"
$ws.Range("E27").Value = "'0.21%"
$ws.Range("D39").Value = "'0.02078"
$ws.Range("E39").Value = "'5.90%"
$ws.Range("D40").Value = "'0.05051"
$ws.Range("E40").Value = "'-2.03%"
$ws.Range("D41").Value = "'0.007469"
$ws.Range("E41").Value = "'-0.44%"
$ws.Range("D42").Value = "'0.01011"
$ws.Range("E42").Value = "'-0.28%"
$ws.Range("E43").Value = "'-0.42%"
$ws.Range("D44").Value = "'0.002131"
$ws.Range("E44").Value = "'1.16%"
$ws.Range("D45").Value = "'0.008938"
$ws.Range("E45").Value = "'-9.71%"
$ws.Range("D46").Value = "'0.00006177"
$ws.Range("E46").Value = "'-0.16%"
$ws.Range("E47").Value = "'0.18%"
$ws.Range("D49").Value = "'0.001600"
$ws.Range("E49").Value = "'28.15%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.18%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.18%"
